$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" (Strike count) values; regenerated with new data.
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 1
